# Update odds figures on the "Jogos da Semana" worksheet to reflect the
# latest FlashScore snapshot (odds and correct-score prices refreshed for
# several fixtures on the 2025-04-15 sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 2.65
$ws.Range("I2").Value = 3.9
$ws.Range("J2").Value = 1.2
$ws.Range("K2").Value = 4.33
$ws.Range("O2").Value = 1.22
$ws.Range("P2").Value = 1.88
$ws.Range("Q2").Value = 1.93
$ws.Range("U2").Value = 9
$ws.Range("AH2").Value = 41
$ws.Range("G3").Value = 2.15
$ws.Range("H3").Value = 3.2
$ws.Range("N3").Value = 2.3
$ws.Range("O3").Value = 1.53
$ws.Range("AF3").Value = 17
$ws.Range("G4").Value = 2.55
$ws.Range("I4").Value = 2.88
$ws.Range("J4").Value = 1.08
$ws.Range("K4").Value = 7.5
$ws.Range("N4").Value = 2.35
$ws.Range("O4").Value = 1.57
$ws.Range("T4").Value = 7
$ws.Range("U4").Value = 11
$ws.Range("V4").Value = 10
$ws.Range("W4").Value = 23
$ws.Range("X4").Value = 23
$ws.Range("Z4").Value = 7.5
$ws.Range("AE4").Value = 7.5
$ws.Range("AF4").Value = 13
$ws.Range("AH4").Value = 29
$ws.Range("AI4").Value = 26
$ws.Range("G8").Value = 1.91
$ws.Range("H8").Value = 3
$ws.Range("J8").Value = 1.11
$ws.Range("K8").Value = 6.5
$ws.Range("N8").Value = 2.6
$ws.Range("O8").Value = 1.48
$ws.Range("P8").Value = 1.57
$ws.Range("Q8").Value = 2.25
$ws.Range("AC8").Value = 81
$ws.Range("G9").Value = 2.7
$ws.Range("H9").Value = 2.82
$ws.Range("I9").Value = 2.75
$ws.Range("L9").Value = 1.5
$ws.Range("M9").Value = 2.5
$ws.Range("N9").Value = 2.6
$ws.Range("O9").Value = 1.48
$ws.Range("P9").Value = 1.57
$ws.Range("Q9").Value = 2.25
$ws.Range("R9").Value = 2.1
$ws.Range("S9").Value = 1.67
$ws.Range("T9").Value = 6.5
$ws.Range("U9").Value = 12
$ws.Range("W9").Value = 29
$ws.Range("AC9").Value = 67
$ws.Range("AE9").Value = 7
$ws.Range("AF9").Value = 12
$ws.Range("AG9").Value = 12
$ws.Range("AH9").Value = 29
$ws.Range("AI9").Value = 29
$ws.Range("AJ9").Value = 41
$ws.Range("G10").Value = 1.91
$ws.Range("J10").Value = 1.11
$ws.Range("K10").Value = 6.5
$ws.Range("L10").Value = 1.5
$ws.Range("M10").Value = 2.5
$ws.Range("P10").Value = 1.57
$ws.Range("Q10").Value = 2.25
$ws.Range("U10").Value = 8
$ws.Range("AG10").Value = 17
$ws.Range("G11").Value = 1.7
$ws.Range("K11").Value = 8
$ws.Range("L11").Value = 1.44
$ws.Range("M11").Value = 2.63
$ws.Range("N11").Value = 2.35
$ws.Range("O11").Value = 1.57
$ws.Range("U11").Value = 7
$ws.Range("AA11").Value = 7
$ws.Range("G12").Value = 1.39
$ws.Range("H12").Value = 4.2
$ws.Range("J12").Value = 1.04
$ws.Range("K12").Value = 13
$ws.Range("L12").Value = 1.2
$ws.Range("M12").Value = 4.33
$ws.Range("N12").Value = 1.7
$ws.Range("O12").Value = 2.1
$ws.Range("P12").Value = 1.33
$ws.Range("Q12").Value = 3.25
$ws.Range("R12").Value = 1.83
$ws.Range("S12").Value = 1.83
$ws.Range("T12").Value = 7.5
$ws.Range("Z12").Value = 13
$ws.Range("AA12").Value = 8.5
$ws.Range("AD12").Value = 301
$ws.Range("AE12").Value = 19
$ws.Range("I13").Value = 1.92
$ws.Range("H14").Value = 3.25
$ws.Range("I14").Value = 2.55
$ws.Range("J14").Value = 1.06
$ws.Range("K14").Value = 10
$ws.Range("P14").Value = 1.4
$ws.Range("Q14").Value = 2.75
$ws.Range("R14").Value = 1.73
$ws.Range("S14").Value = 2
$ws.Range("T14").Value = 9
$ws.Range("AB14").Value = 13
$ws.Range("AC14").Value = 41
$ws.Range("AD14").Value = 201
$ws.Range("AE14").Value = 9
$ws.Range("G18").Value = 1.86
$ws.Range("L18").Value = 1.33
$ws.Range("M18").Value = 3.25
$ws.Range("O18").Value = 1.72
$ws.Range("J19").Value = 1.05
$ws.Range("K19").Value = 11
$ws.Range("N19").Value = 1.93
$ws.Range("O19").Value = 1.93
$ws.Range("G20").Value = 2.55
$ws.Range("I20").Value = 2.62
$ws.Range("Z20").Value = 11
$ws.Range("AE20").Value = 9
$ws.Range("AH20").Value = 26
